# Add a new "Jun-2017" monthly column ahead of the existing month columns.
# The report's monthly columns (J:U) each hold one calendar month of data,
# newest first. A new month rolled over, so every existing month shifts one
# column to the right and the oldest month (the last column, U) falls off
# the right edge. The header row gets a brand-new "Jun-2017" label in J1 and
# the data rows get a fresh 0 in column J (no activity recorded yet for the
# new month).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): shift text labels J1:U1 right by one column ----
# Plain Value assignment of a string like "May-2017" gets auto-coerced by
# Excel into a date serial, which would corrupt the header's text type and
# style. Copy/PasteSpecial(xlPasteValues) moves the already-typed string
# value verbatim, without re-parsing the display text, so it stays text.
$headerCols = @("J","K","L","M","N","O","P","Q","R","S","T","U")
for ($i = $headerCols.Length - 1; $i -ge 1; $i--) {
    $srcCell = $headerCols[$i - 1] + "1"
    $dstCell = $headerCols[$i] + "1"
    $ws.Range($srcCell).Copy()
    $ws.Range($dstCell).PasteSpecial(-4163)  # xlPasteValues
}

# New J1 needs literal new text that doesn't exist anywhere else in the
# sheet yet. Temporarily mark the cell as Text so the assignment isn't
# reinterpreted as a date, then restore the original (General) number
# format by pulling it back from a neighboring header cell (I1) so the
# cell's style matches the rest of the header row exactly.
$ws.Range("J1").NumberFormat = "@"
$ws.Range("J1").Value = "Jun-2017"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# ---- Data rows (rows 2:73): shift the monthly counts right by one column ----
for ($r = 2; $r -le 73; $r++) {
    $vals = @()
    for ($col = 10; $col -le 21; $col++) {
        $v = $ws.Cells.Item($r, $col).Value2
        if ($null -eq $v) { $v = 0 }
        $vals += $v
    }
    for ($i = 21; $i -ge 11; $i--) {
        $ws.Cells.Item($r, $i).Value2 = $vals[$i - 11]
    }
    $ws.Cells.Item($r, 10).Value2 = 0
}
